# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to match the refreshed scrape output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 15054
$ws1.Range("F3").Value  = 19178
$ws1.Range("F7").Value  = 223
$ws1.Range("F14").Value = 172
$ws1.Range("F17").Value = 1484
$ws1.Range("F22").Value = 8006
$ws1.Range("F24").Value = 34
$ws1.Range("F25").Value = 3
$ws1.Range("F29").Value = 6077
$ws1.Range("F31").Value = 73
$ws1.Range("F32").Value = 173
$ws1.Range("F35").Value = 5475
$ws1.Range("F36").Value = 765

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 15054
$ws4.Range("F3").Value  = 19178
$ws4.Range("F7").Value  = 223
$ws4.Range("F14").Value = 172
$ws4.Range("F17").Value = 1484
$ws4.Range("F23").Value = 8006
$ws4.Range("F25").Value = 34
$ws4.Range("F26").Value = 3
$ws4.Range("F32").Value = 6077
$ws4.Range("F34").Value = 73
$ws4.Range("F35").Value = 173
$ws4.Range("F38").Value = 5475
$ws4.Range("F39").Value = 765
